$d = $word.ActiveDocument

# The last paragraph currently holds the "_GoBack" bookmark right after its
# run. Two new bulleted (Listenabsatz / numId 1) paragraphs need to be
# appended after it, and the bookmark needs to end up after the very last
# of those new paragraphs.

$last = $d.Paragraphs.Last
$r = $last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$p2 = $d.Paragraphs.Last
$p2.Range.Font.Bold = 0
$p2.Range.Text = "Maybe add short fills/premade juice to the DB as well?"

$r2 = $p2.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()

$p3 = $d.Paragraphs.Last
$p3.Range.Font.Bold = 0
$p3.Range.Text = "Auto suggestions for recipes based on the stuff a user likes and mixes often, as well as based on their stash (“We see you’ve been mixing a lot of recipes including strawberry! Maybe you want to try THIS recipe”). You could also go a little extra step and only suggest recipes that the user can mix right away, based on their current stash"

Write-Output $d.Paragraphs.Count
